$d = $word.ActiveDocument
$para = $d.Paragraphs(1)
$para.Range.Font.Color = 255
